$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb1"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.46510533333333
$ws.Range("H2").Value = 31.395316
$ws.Range("I2").Value = 0.5554075997074865
$ws.Range("J2").Value = 0.5554075997074865
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.885873333333334
$ws.Range("N2").Value = 8.657620000000001
$ws.Range("O2").Value = 0.3070415651026022
$ws.Range("P2").Value = 0.3070415651026022
$ws.Range("Q2").Value = 30.20096841199112
$ws.Range("R2").Value = 271.80871570792
$ws.Range("S2").Value = 0.1705332186840663
$ws.Range("T2").Value = 0.1705332186840663

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb1"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.46510533333333
$ws.Range("H3").Value = 31.395316
$ws.Range("I3").Value = 0.5554075997074865
$ws.Range("J3").Value = 0.5554075997074865
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.165953666666667
$ws.Range("N3").Value = 9.497861
$ws.Range("O3").Value = 0.3368406220840099
$ws.Range("P3").Value = 0.3368406220840099
$ws.Range("Q3").Value = 33.13203860211956
$ws.Range("R3").Value = 298.188347419076
$ws.Range("S3").Value = 0.1870838413956565
$ws.Range("T3").Value = 0.1870838413956565

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efnb1"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.46510533333333
$ws.Range("H4").Value = 31.395316
$ws.Range("I4").Value = 0.5554075997074865
$ws.Range("J4").Value = 0.5554075997074865
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.327024333333334
$ws.Range("N4").Value = 9.981073
$ws.Range("O4").Value = 0.3539776838580724
$ws.Range("P4").Value = 0.3539776838580724
$ws.Range("Q4").Value = 34.81766009489645
$ws.Range("R4").Value = 313.358940854068
$ws.Range("S4").Value = 0.1966018957416275
$ws.Range("T4").Value = 0.1966018957416275

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Efnb1"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 10.46510533333333
$ws.Range("H5").Value = 31.395316
$ws.Range("I5").Value = 0.5554075997074865
$ws.Range("J5").Value = 0.5554075997074865
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.020115
$ws.Range("N5").Value = 0.060345
$ws.Range("O5").Value = 0.002140128955315263
$ws.Range("P5").Value = 0.002140128955315263
$ws.Range("Q5").Value = 0.21050559378
$ws.Range("R5").Value = 1.89455034402
$ws.Range("S5").Value = 0.001188643886136141
$ws.Range("T5").Value = 0.001188643886136141

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efnb1"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.754308333333334
$ws.Range("H6").Value = 17.262925
$ws.Range("I6").Value = 0.3053945925621632
$ws.Range("J6").Value = 0.3053945925621632
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.885873333333334
$ws.Range("N6").Value = 8.657620000000001
$ws.Range("O6").Value = 0.3070415651026022
$ws.Range("P6").Value = 0.3070415651026022
$ws.Range("Q6").Value = 16.60620497094445
$ws.Range("R6").Value = 149.4558447385
$ws.Range("S6").Value = 0.0937688336741581
$ws.Range("T6").Value = 0.0937688336741581

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efnb1"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.754308333333334
$ws.Range("H7").Value = 17.262925
$ws.Range("I7").Value = 0.3053945925621632
$ws.Range("J7").Value = 0.3053945925621632
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.165953666666667
$ws.Range("N7").Value = 9.497861
$ws.Range("O7").Value = 0.3368406220840099
$ws.Range("P7").Value = 0.3368406220840099
$ws.Range("Q7").Value = 18.21787356704722
$ws.Range("R7").Value = 163.960862103425
$ws.Range("S7").Value = 0.1028693045397318
$ws.Range("T7").Value = 0.1028693045397318

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Efnb1"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.754308333333334
$ws.Range("H8").Value = 17.262925
$ws.Range("I8").Value = 0.3053945925621632
$ws.Range("J8").Value = 0.3053945925621632
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.327024333333334
$ws.Range("N8").Value = 9.981073
$ws.Range("O8").Value = 0.3539776838580724
$ws.Range("P8").Value = 0.3539776838580724
$ws.Range("Q8").Value = 19.14472384650278
$ws.Range("R8").Value = 172.302514618525
$ws.Range("S8").Value = 0.1081028705379342
$ws.Range("T8").Value = 0.1081028705379342

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Efnb1"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.754308333333334
$ws.Range("H9").Value = 17.262925
$ws.Range("I9").Value = 0.3053945925621632
$ws.Range("J9").Value = 0.3053945925621632
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.020115
$ws.Range("N9").Value = 0.060345
$ws.Range("O9").Value = 0.002140128955315263
$ws.Range("P9").Value = 0.002140128955315263
$ws.Range("Q9").Value = 0.115747912125
$ws.Range("R9").Value = 1.041731209125
$ws.Range("S9").Value = 0.0006535838103389926
$ws.Range("T9").Value = 0.0006535838103389927

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Efnb1"
$ws.Range("C10").Value = "Erbb2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.146766
$ws.Range("H10").Value = 6.440298
$ws.Range("I10").Value = 0.1139338891693565
$ws.Range("J10").Value = 0.1139338891693565
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.885873333333334
$ws.Range("N10").Value = 8.657620000000001
$ws.Range("O10").Value = 0.3070415651026022
$ws.Range("P10").Value = 0.3070415651026022
$ws.Range("Q10").Value = 6.195294752306668
$ws.Range("R10").Value = 55.75765277076001
$ws.Range("S10").Value = 0.03498243964878565
$ws.Range("T10").Value = 0.03498243964878565

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Efnb1"
$ws.Range("C11").Value = "Erbb2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.146766
$ws.Range("H11").Value = 6.440298
$ws.Range("I11").Value = 0.1139338891693565
$ws.Range("J11").Value = 0.1139338891693565
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.165953666666667
$ws.Range("N11").Value = 9.497861
$ws.Range("O11").Value = 0.3368406220840099
$ws.Range("P11").Value = 0.3368406220840099
$ws.Range("Q11").Value = 6.796561689175333
$ws.Range("R11").Value = 61.16905520257801
$ws.Range("S11").Value = 0.0383775621042567
$ws.Range("T11").Value = 0.0383775621042567

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Efnb1"
$ws.Range("C12").Value = "Erbb2"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.146766
$ws.Range("H12").Value = 6.440298
$ws.Range("I12").Value = 0.1139338891693565
$ws.Range("J12").Value = 0.1139338891693565
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 3.327024333333334
$ws.Range("N12").Value = 9.981073
$ws.Range("O12").Value = 0.3539776838580724
$ws.Range("P12").Value = 0.3539776838580724
$ws.Range("Q12").Value = 7.142342719972667
$ws.Range("R12").Value = 64.281084479754
$ws.Range("S12").Value = 0.04033005420111115
$ws.Range("T12").Value = 0.04033005420111115

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Efnb1"
$ws.Range("C13").Value = "Erbb2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.146766
$ws.Range("H13").Value = 6.440298
$ws.Range("I13").Value = 0.1139338891693565
$ws.Range("J13").Value = 0.1139338891693565
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.020115
$ws.Range("N13").Value = 0.060345
$ws.Range("O13").Value = 0.002140128955315263
$ws.Range("P13").Value = 0.002140128955315263
$ws.Range("Q13").Value = 0.04318219809
$ws.Range("R13").Value = 0.3886397828100001
$ws.Range("S13").Value = 0.0002438332152030199
$ws.Range("T13").Value = 0.00024383321520302

# Row 14
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Efnb1"
$ws.Range("C14").Value = "Erbb2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.4760280000000001
$ws.Range("H14").Value = 1.428084
$ws.Range("I14").Value = 0.02526391856099382
$ws.Range("J14").Value = 0.02526391856099382
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.885873333333334
$ws.Range("N14").Value = 8.657620000000001
$ws.Range("O14").Value = 0.3070415651026022
$ws.Range("P14").Value = 0.3070415651026022
$ws.Range("Q14").Value = 1.373756511120001
$ws.Range("R14").Value = 12.36380860008
$ws.Range("S14").Value = 0.007757073095592225
$ws.Range("T14").Value = 0.007757073095592225

# Row 15
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Efnb1"
$ws.Range("C15").Value = "Erbb2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.4760280000000001
$ws.Range("H15").Value = 1.428084
$ws.Range("I15").Value = 0.02526391856099382
$ws.Range("J15").Value = 0.02526391856099382
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 3.165953666666667
$ws.Range("N15").Value = 9.497861
$ws.Range("O15").Value = 0.3368406220840099
$ws.Range("P15").Value = 0.3368406220840099
$ws.Range("Q15").Value = 1.507082592036
$ws.Range("R15").Value = 13.563743328324
$ws.Range("S15").Value = 0.008509914044364924
$ws.Range("T15").Value = 0.008509914044364924

# Row 16
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Efnb1"
$ws.Range("C16").Value = "Erbb2"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.4760280000000001
$ws.Range("H16").Value = 1.428084
$ws.Range("I16").Value = 0.02526391856099382
$ws.Range("J16").Value = 0.02526391856099382
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 3.327024333333334
$ws.Range("N16").Value = 9.981073
$ws.Range("O16").Value = 0.3539776838580724
$ws.Range("P16").Value = 0.3539776838580724
$ws.Range("Q16").Value = 1.583756739348
$ws.Range("R16").Value = 14.253810654132
$ws.Range("S16").Value = 0.008942863377399558
$ws.Range("T16").Value = 0.008942863377399558

# Row 17
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Efnb1"
$ws.Range("C17").Value = "Erbb2"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.4760280000000001
$ws.Range("H17").Value = 1.428084
$ws.Range("I17").Value = 0.02526391856099382
$ws.Range("J17").Value = 0.02526391856099382
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.020115
$ws.Range("N17").Value = 0.060345
$ws.Range("O17").Value = 0.002140128955315263
$ws.Range("P17").Value = 0.002140128955315263
$ws.Range("Q17").Value = 0.009575303220000001
$ws.Range("R17").Value = 0.08617772898
$ws.Range("S17").Value = 0.00005406804363710958
$ws.Range("T17").Value = 0.0000540680436371096

